$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dates = @(
    "Wed Nov 01 15:37:06 EDT 2023",
    "Wed Nov 01 15:37:16 EDT 2023",
    "Wed Nov 01 15:37:26 EDT 2023",
    "Wed Nov 01 15:37:36 EDT 2023",
    "Wed Nov 01 15:37:45 EDT 2023",
    "Wed Nov 01 15:37:55 EDT 2023",
    "Wed Nov 01 15:38:05 EDT 2023",
    "Wed Nov 01 15:38:15 EDT 2023",
    "Wed Nov 01 15:38:25 EDT 2023",
    "Wed Nov 01 15:38:35 EDT 2023",
    "Wed Nov 01 15:38:45 EDT 2023",
    "Wed Nov 01 15:38:56 EDT 2023"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}
